# expansão das análises automáticas
# Adds three new computed columns (apoio_medio, contribuicoes, media_contribuicoes)
# and rescales the two existing percentage columns (taxa de sucesso / particip)
# from fractional (0-1) to percentage-point (0-100) numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column headers (L, M, N) ------------------------------------------------
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- Rescale existing percentage columns E (particip) and F (taxa_sucesso) -------
# Previously stored as a 0-1 fraction; now stored as a 0-100 percentage number
# (the cells keep their existing "0.00%" number format).
$ws.Range("E2").Value = 79.32584269662921
$ws.Range("F2").Value = 63.5505193578848

$ws.Range("E3").Value = 20.67415730337078
$ws.Range("F3").Value = 56.88405797101449

$ws.Range("E4").Value = 75.88555858310627
$ws.Range("F4").Value = 93.35727109515261

$ws.Range("E5").Value = 24.11444141689373
$ws.Range("F5").Value = 96.89265536723164

$ws.Range("E6").Value = 93.56725146198829
$ws.Range("F6").Value = 22.03125

$ws.Range("E7").Value = 6.432748538011696
$ws.Range("F7").Value = 25

# --- New column data (L = apoio_medio, M = contribuicoes, N = media_contribuicoes)
$ws.Range("L2").Value = 91.32270166935785
$ws.Range("M2").Value = 225451
$ws.Range("N2").Value = 334.9940564635958

$ws.Range("L3").Value = 91.18905604074934
$ws.Range("M3").Value = 38102
$ws.Range("N3").Value = 242.687898089172

$ws.Range("L4").Value = 87.48757423640602
$ws.Range("M4").Value = 147045
$ws.Range("N4").Value = 141.3894230769231

$ws.Range("L5").Value = 97.12764057113438
$ws.Range("M5").Value = 56601
$ws.Range("N5").Value = 165.0174927113703

$ws.Range("L6").Value = 18.13544486308301
$ws.Range("M6").Value = 1980
$ws.Range("N6").Value = 14.04255319148936

$ws.Range("L7").Value = 31.92446020123002
$ws.Range("M7").Value = 228
$ws.Range("N7").Value = 20.72727272727273
